$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# -----------------------------------------------------------------
# 1. Insert a new blank row at row 13. This pushes the "Camel
#    Definitions" section header (and everything below it) down by
#    one row and grows the used range from B29 to B30. Inserting at
#    row 13 (rather than row 12) makes the new row pick up the
#    formatting of row 12 -- which is a lone styled, empty A cell
#    with no B cell at all -- exactly matching the target's new
#    row 13 (only "A13" with style 1, no "B13").
# -----------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# -----------------------------------------------------------------
# 2. Swap the "Eclipse Workspace for Camel + ActiveMQ Examples" row
#    (row 10) with the "Github Repository for Camel + ActiveMQ
#    Examples" row (row 11): same two rows, values/styles/hyperlink
#    traded places.
# -----------------------------------------------------------------

# -- Column A: plain text values, nothing else to carry along. --
$a10 = $ws.Range("A10").Value2
$a11 = $ws.Range("A11").Value2
$ws.Range("A10").Value2 = $a11
$ws.Range("A11").Value2 = $a10

# -- Column B: swap the whole cell (value + style) through a scratch
#    cell with Copy so the existing "Hyperlink" style record (used by
#    B11) relocates to B10 instead of a new style being minted. --
$scratch = $ws.Range("D100")
$styleScratch = $ws.Range("D101")

$ws.Range("B10").Copy($scratch)
$ws.Range("B11").Copy($ws.Range("B10"))
$scratch.Copy($ws.Range("B11"))
$scratch.Clear()

# Stash away a copy of B10's now-correct formatting (style index 7)
# before Hyperlinks.Add (below) gets a chance to overwrite it with its
# own hyperlink style variant.
$ws.Range("B10").Copy($styleScratch)

# -- Hyperlink bookkeeping: the Hyperlink object is anchored to a
#    range and does not follow a value/style copy on its own, so move
#    it explicitly from B11 to B10. --
$oldHyperlink = $null
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$11') {
        $oldHyperlink = $hl
    }
}
$hyperlinkAddress = $oldHyperlink.Address
$hyperlinkText = $ws.Range("B10").Value2
$oldHyperlink.Delete()
$ws.Hyperlinks.Add($ws.Range("B10"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $hyperlinkText) | Out-Null

# Restore B10's exact original formatting (Hyperlinks.Add forces its
# own style variant onto the cell); keep the new value/hyperlink.
$styleScratch.Copy()
$ws.Range("B10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$styleScratch.Clear()

# -----------------------------------------------------------------
# 3. Update the active selection to A11, matching the saved view
#    state in the target workbook.
# -----------------------------------------------------------------
$ws.Range("A11").Select()
